$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update driving parameter D14 (switching frequency) 40000 -> 50000
$ws.Range("D14").Value = 50000

# Update "Enter Switching Energy Here for 150C" block (A22:D22) with the
# Ron=15 switching-energy values
$ws.Range("A22").Value = 0.0000408
$ws.Range("B22").Value = 0.000005149
$ws.Range("C22").Value = -0.0000042
$ws.Range("D22").Value = 0.00000528

# Pon/Poff formulas no longer multiply by the switching frequency column C14
$ws.Range("F3").Formula = "=(A19+C19)*D14/(PI()*2)"
$ws.Range("G3").Formula = "=(B19+D19)*D14/PI()/2"
$ws.Range("F8").Formula = "=(A22+C22)*D14/PI()/2"
$ws.Range("G8").Formula = "=(B22+D22)*D14/PI()/2"

# D19 right border was thin; make it a proper medium border, matching the
# rest of the closed table boxes on the sheet
$ws.Range("D19").Borders.Item(10).LineStyle = 1
$ws.Range("D19").Borders.Item(10).Weight = -4138

# Restore the active selection to L16
$ws.Range("L16").Select()
